$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" section:
# "LOQ4205: Sistemas Produtivos II (Requisito fraco)"
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "LOQ4205: Sistemas Produtivos II (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'LOQ4205' requisito paragraph."
}

$anchorIndex = $anchor.Paragraphs.Item(1).Index

# Immediately after it come three paragraphs that must be removed:
#   1) a blank spacer paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "(c) 2020 ... Creative Commons Attribution" footer line
$firstToRemove = $d.Paragraphs.Item($anchorIndex + 1)
$lastToRemove  = $d.Paragraphs.Item($anchorIndex + 3)

$deleteRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)

$deleteRange.Delete()
